$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kml data")

$ws.Range("R2").Value = "BRS"
$ws.Range("S2").Value = "JWA"
$ws.Range("R3").Value = "CMC"
$ws.Range("S3").Value = "JMO"
$ws.Range("R4").Value = "JBB"
$ws.Range("S4").Value = "JPJ"
$ws.Range("Q5").Value = "2026-02-10"
$ws.Range("R5").Value = "BRS"
$ws.Range("S5").Value = "JWA"
$ws.Range("Q6").Value = "2026-02-09"
$ws.Range("R6").Value = "LDE"
$ws.Range("S6").Value = "WWH"
$ws.Range("R7").Value = "CMC"
$ws.Range("S7").Value = "PRT"
$ws.Range("Q8").Value = "2026-02-06"
$ws.Range("R8").Value = "KCO"
$ws.Range("S8").Value = "WWH"
$ws.Range("R9").Value = "AFU"
$ws.Range("S9").Value = "BAP"
$ws.Range("R10").Value = "CJO"
$ws.Range("S10").Value = "CMC"
$ws.Range("Q11").Value = "2026-02-05"
$ws.Range("R11").Value = "AFU"
$ws.Range("S11").Value = "BAP"
$ws.Range("R13").Value = "HAS"
$ws.Range("S13").Value = "WSM"
$ws.Range("R14").Value = "BRS"
$ws.Range("S14").Value = ""
$ws.Range("Q16").Value = "2026-02-04"
$ws.Range("R16").Value = "DDC"
$ws.Range("S16").Value = "RBL"
$ws.Range("Q17").Value = "2026-02-02"
$ws.Range("R17").Value = "CBI"
$ws.Range("S17").Value = "LDE"
$ws.Range("Q18").Value = "2026-01-29"
$ws.Range("R18").Value = "SMB"
$ws.Range("Q19").Value = "2026-01-22"
$ws.Range("R19").Value = "AWP"
$ws.Range("S19").Value = "JWG"
$ws.Range("Q20").Value = "2026-01-21"
$ws.Range("R20").Value = "JPJ"
$ws.Range("S20").Value = "WLS"
$ws.Range("Q21").Value = "2026-01-18"
$ws.Range("S21").Value = "HAS"
$ws.Range("Q22").Value = "2026-01-12"
$ws.Range("R22").Value = "ADO"
$ws.Range("S22").Value = "GLO"
$ws.Range("Q23").Value = "2026-01-05"
$ws.Range("R23").Value = "JWA"
$ws.Range("S23").Value = "MAN"
$ws.Range("Q25").Value = "2025-12-18"
$ws.Range("R25").Value = "CBI"
$ws.Range("Q26").Value = "2025-12-15"
$ws.Range("R26").Value = "CAD"
$ws.Range("S26").Value = "LDE"
$ws.Range("Q27").Value = "2025-12-11"
$ws.Range("R27").Value = "GLO"
$ws.Range("S27").Value = "HMK"
$ws.Range("Q28").Value = "2025-12-09"
$ws.Range("R28").Value = "JMO"
$ws.Range("S28").Value = "RJC"
$ws.Range("Q29").Value = "2025-12-02"
$ws.Range("R29").Value = "BRS"
$ws.Range("S29").Value = "LOV"
$ws.Range("R30").Value = "CAD"
$ws.Range("S30").Value = "TIN"
$ws.Range("Q31").Value = "2025-11-21"
$ws.Range("R31").Value = "CBI"
$ws.Range("S31").Value = "SMB"
$ws.Range("R32").Value = "MAN"
$ws.Range("S32").Value = ""
$ws.Range("Q33").Value = "2025-11-19"
$ws.Range("R33").Value = "BRS"
$ws.Range("S33").Value = "FRU"
$ws.Range("Q34").Value = "2025-11-14"
$ws.Range("R34").Value = "LAO"
$ws.Range("S34").Value = "RJC"
$ws.Range("Q35").Value = "2025-11-06"
$ws.Range("R35").Value = "GLO"
$ws.Range("S35").Value = ""
$ws.Range("Q36").Value = "2025-11-05"
$ws.Range("R36").Value = "AWP"
$ws.Range("S36").Value = "BNE"
$ws.Range("Q37").Value = "2025-11-04"
$ws.Range("R37").Value = "DAR"
$ws.Range("S37").Value = "RTH"
$ws.Range("Q40").Value = "2025-10-26"
$ws.Range("R40").Value = "DDC"
$ws.Range("Q41").Value = "2025-10-15"
$ws.Range("R41").Value = "GLO"
$ws.Range("S41").Value = ""
$ws.Range("Q42").Value = "2025-10-13"
$ws.Range("R42").Value = "BAP"
$ws.Range("S42").Value = "FRU"
$ws.Range("R48").Value = "MBE"
$ws.Range("S48").Value = "RJC"
$ws.Range("R66").Value = "JWG"
$ws.Range("S66").Value = "SMB"
$ws.Range("R68").Value = "RJC"
$ws.Range("R69").Value = "MAN"
$ws.Range("S80").Value = "WLS"
$ws.Range("S81").Value = ""
$ws.Range("R84").Value = "ASM"
$ws.Range("S84").Value = "CMC"
$ws.Range("R85").Value = "AWP"
$ws.Range("S85").Value = ""
$ws.Range("R105").Value = "LOV"
$ws.Range("S105").Value = "RTH"
$ws.Range("R106").Value = "DDC"
$ws.Range("S106").Value = ""
